$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "/RME" from the steel (S) line of the industrial mapping description.
$newText = "6% MUR/LWAL+CDN/H:1`n14% CR/LFM+CDN/H:2`n31% CR+PC/LFM+CDN/H:1`n33% S/LFM+CDN/H:1`n13% S+SL/LFM+CDN/H:1`n3% S/LFBR+CDN/H:1"
$ws.Range("B2").Value = $newText

# Wrap the long multi-line description and grow the row to show all 6 lines.
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 320

[void]$ws.Range("B9").Select()
